$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08805033333333334
$ws.Range("H2").Value = 0.264151
$ws.Range("I2").Value = 0.005589762818257384
$ws.Range("J2").Value = 0.005589762818257385
$ws.Range("M2").Value = 8.252454666666667
$ws.Range("N2").Value = 24.757364
$ws.Range("O2").Value = 0.05349680956196952
$ws.Range("P2").Value = 0.05349680956196953
$ws.Range("Q2").Value = 0.7266313842182223
$ws.Range("R2").Value = 6.539682457964001
$ws.Range("S2").Value = 0.0002990344769848933
$ws.Range("T2").Value = 0.0002990344769848934
$ws.Range("G3").Value = 0.08805033333333334
$ws.Range("H3").Value = 0.264151
$ws.Range("I3").Value = 0.005589762818257384
$ws.Range("J3").Value = 0.005589762818257385
$ws.Range("O3").Value = 0.5638948237978928
$ws.Range("P3").Value = 0.5638948237978929
$ws.Range("Q3").Value = 7.659217058451223
$ws.Range("R3").Value = 68.93295352606101
$ws.Range("S3").Value = 0.003152038319473261
$ws.Range("T3").Value = 0.003152038319473261
$ws.Range("G4").Value = 0.08805033333333334
$ws.Range("H4").Value = 0.264151
$ws.Range("I4").Value = 0.005589762818257384
$ws.Range("J4").Value = 0.005589762818257385
$ws.Range("M4").Value = 57.81408433333333
$ws.Range("N4").Value = 173.442253
$ws.Range("O4").Value = 0.3747817085348802
$ws.Range("P4").Value = 0.3747817085348802
$ws.Range("Q4").Value = 5.090549396911445
$ws.Range("R4").Value = 45.814944572203
$ws.Range("S4").Value = 0.002094940859331249
$ws.Range("T4").Value = 0.00209494085933125
$ws.Range("G5").Value = 0.08805033333333334
$ws.Range("H5").Value = 0.264151
$ws.Range("I5").Value = 0.005589762818257384
$ws.Range("J5").Value = 0.005589762818257385
$ws.Range("M5").Value = 1.207345666666667
$ws.Range("N5").Value = 3.622037
$ws.Range("O5").Value = 0.007826658105257385
$ws.Range("P5").Value = 0.007826658105257386
$ws.Range("Q5").Value = 0.1063071883985556
$ws.Range("R5").Value = 0.9567646955870002
$ws.Range("S5").Value = 0.00004374916246798052
$ws.Range("T5").Value = 0.00004374916246798053
$ws.Range("I6").Value = 0.9470512964761942
$ws.Range("J6").Value = 0.9470512964761943
$ws.Range("M6").Value = 8.252454666666667
$ws.Range("N6").Value = 24.757364
$ws.Range("O6").Value = 0.05349680956196952
$ws.Range("P6").Value = 0.05349680956196953
$ws.Range("Q6").Value = 123.1102672615174
$ws.Range("R6").Value = 1107.992405353656
$ws.Range("S6").Value = 0.05066422285300329
$ws.Range("T6").Value = 0.05066422285300331
$ws.Range("I7").Value = 0.9470512964761942
$ws.Range("J7").Value = 0.9470512964761943
$ws.Range("O7").Value = 0.5638948237978928
$ws.Range("P7").Value = 0.5638948237978929
$ws.Range("Q7").Value = 1297.6707028618
$ws.Range("R7").Value = 11679.0363257562
$ws.Range("S7").Value = 0.5340373239540095
$ws.Range("T7").Value = 0.5340373239540096
$ws.Range("I8").Value = 0.9470512964761942
$ws.Range("J8").Value = 0.9470512964761943
$ws.Range("M8").Value = 57.81408433333333
$ws.Range("N8").Value = 173.442253
$ws.Range("O8").Value = 0.3747817085348802
$ws.Range("P8").Value = 0.3747817085348802
$ws.Range("Q8").Value = 862.4715507381848
$ws.Range("R8").Value = 7762.243956643662
$ws.Range("S8").Value = 0.3549375029635214
$ws.Range("T8").Value = 0.3549375029635215
$ws.Range("I9").Value = 0.9470512964761942
$ws.Range("J9").Value = 0.9470512964761943
$ws.Range("M9").Value = 1.207345666666667
$ws.Range("N9").Value = 3.622037
$ws.Range("O9").Value = 0.007826658105257385
$ws.Range("P9").Value = 0.007826658105257386
$ws.Range("Q9").Value = 18.01120438755533
$ws.Range("R9").Value = 162.100839487998
$ws.Range("S9").Value = 0.00741224670565992
$ws.Range("T9").Value = 0.007412246705659923
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7460013333333334
$ws.Range("H10").Value = 2.238004
$ws.Range("I10").Value = 0.04735894070554834
$ws.Range("J10").Value = 0.04735894070554835
$ws.Range("M10").Value = 8.252454666666667
$ws.Range("N10").Value = 24.757364
$ws.Range("O10").Value = 0.05349680956196952
$ws.Range("P10").Value = 0.05349680956196953
$ws.Range("Q10").Value = 6.156342184606223
$ws.Range("R10").Value = 55.40707966145601
$ws.Range("S10").Value = 0.002533552231981326
$ws.Range("T10").Value = 0.002533552231981327
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.7460013333333334
$ws.Range("H11").Value = 2.238004
$ws.Range("I11").Value = 0.04735894070554834
$ws.Range("J11").Value = 0.04735894070554835
$ws.Range("O11").Value = 0.5638948237978928
$ws.Range("P11").Value = 0.5638948237978929
$ws.Range("Q11").Value = 64.89227151773824
$ws.Range("R11").Value = 584.0304436596441
$ws.Range("S11").Value = 0.02670546152441004
$ws.Range("T11").Value = 0.02670546152441005
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.7460013333333334
$ws.Range("H12").Value = 2.238004
$ws.Range("I12").Value = 0.04735894070554834
$ws.Range("J12").Value = 0.04735894070554835
$ws.Range("M12").Value = 57.81408433333333
$ws.Range("N12").Value = 173.442253
$ws.Range("O12").Value = 0.3747817085348802
$ws.Range("P12").Value = 0.3747817085348802
$ws.Range("Q12").Value = 43.12938399811245
$ws.Range("R12").Value = 388.164455983012
$ws.Range("S12").Value = 0.01774926471202749
$ws.Range("T12").Value = 0.0177492647120275
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.7460013333333334
$ws.Range("H13").Value = 2.238004
$ws.Range("I13").Value = 0.04735894070554834
$ws.Range("J13").Value = 0.04735894070554835
$ws.Range("M13").Value = 1.207345666666667
$ws.Range("N13").Value = 3.622037
$ws.Range("O13").Value = 0.007826658105257385
$ws.Range("P13").Value = 0.007826658105257386
$ws.Range("Q13").Value = 0.9006814771275556
$ws.Range("R13").Value = 8.106133294148
$ws.Range("S13").Value = 0.0003706622371294838
$ws.Range("T13").Value = 0.0003706622371294834
